# QAPF map area 3c
# Adds a new "Area 3c" data block (columns AZ:BG) mirroring the existing
# "Area 3a" block (columns AQ:AX), labels the AQ column header "Area 3b",
# adds a new "area 3c" search-radius-gain column (L), and updates the
# sheet selection / column width accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New shared-string-backed labels.
#    Insertion order matters for shared-string table ordering, so we
#    write AZ18 ("Area 3c") first, then AQ18 ("Area 3b"), then L18
#    ("area 3c") - matching the order the strings appear in the target
#    workbook's sst.
# ---------------------------------------------------------------------
$ws.Range("AZ18").Value2 = "Area 3c"
$ws.Range("AQ18").Value2 = "Area 3b"
$ws.Range("L18").Value2 = "area 3c"

# ---------------------------------------------------------------------
# 2. Formatting: copy the existing "Area" title style (bold, fill,
#    left border only - same look as the other area-group headers in
#    row 2 / row 18) onto the two new title cells.
# ---------------------------------------------------------------------
$ws.Range("AZ2").Copy() | Out-Null
$ws.Range("AQ18,AZ18").PasteSpecial(-4122) | Out-Null

# Column sub-headers (search_radius,#PCs,VAR(OK),MSPE,S_nugget,
# VAR(TOTAL),VAR(DATA)) for the new Area 3c block - copy the format of
# the matching existing sub-header row.
$ws.Range("AI18:AO18").Copy() | Out-Null
$ws.Range("BA18").PasteSpecial(-4122) | Out-Null

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Sub-header text for the new Area 3c block (BA18:BG18), matching
#    the existing Area 3a block (AR18:AX18) labels.
# ---------------------------------------------------------------------
$ws.Range("BA18").Value2 = "search_radius"
$ws.Range("BB18").Value2 = "#PCs"
$ws.Range("BC18").Value2 = "VAR(OK)"
$ws.Range("BD18").Value2 = "MSPE"
$ws.Range("BE18").Value2 = "S_nugget"
$ws.Range("BF18").Value2 = "VAR(TOTAL)"
$ws.Range("BG18").Value2 = "VAR(DATA)"

# ---------------------------------------------------------------------
# 4. Data rows 19-27 for the new Area 3c block (columns AZ:BG).
#    (BF - VAR(TOTAL) - is left blank, matching the Area 3a block.)
# ---------------------------------------------------------------------
$az = @(0,1,2,3,4,5,6,7,8)
$ba = @(20000,20000,20000,20000,20000,20000,20000,20000,20000)
$bb = @(1,2,3,4,5,6,7,8,9)
$bc = @(0.2112798674733572,0.78896175297514703,0.85624013185044456,0.90795121919934418,0.92721678724985468,0.93248005872858164,0.9374691948281535,0.94016820150631564,0.94065254115938446)
$bd = @(3.6981760732177049,3.3784626087081819,3.390308759630241,3.3820575080324442,3.4175163136841169,3.432024431040813,3.435812110388047,3.4360417969815131,3.4367920711214688)
$be = @(1.5298,1.649,1.67299,1.7728900000000001,1.8817900000000001,1.90463,1.92791,1.9418500000000001,1.9454979999999999)
$bg = @(3.5115795310554612,3.5115795310554612,3.5115795310554612,3.5115795310554612,3.5115795310554612,3.5115795310554612,3.5115795310554612,3.5115795310554612,3.5115795310554612)

for ($i = 0; $i -lt 9; $i++) {
    $row = 19 + $i
    $ws.Cells.Item($row, 52).Value2 = $az[$i]
    $ws.Cells.Item($row, 53).Value2 = $ba[$i]
    $ws.Cells.Item($row, 54).Value2 = $bb[$i]
    $ws.Cells.Item($row, 55).Value2 = $bc[$i]
    $ws.Cells.Item($row, 56).Value2 = $bd[$i]
    $ws.Cells.Item($row, 57).Value2 = $be[$i]
    $ws.Cells.Item($row, 59).Value2 = $bg[$i]
}

# Copy the AZ19:AZ27 "index" style (matches AQ19:AQ27 / AH19:AH27).
$ws.Range("AH19").Copy() | Out-Null
$ws.Range("AZ19:AZ27").PasteSpecial(-4122) | Out-Null
$ws.Range("AQ19:AQ27").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 5. Column L ("area 3c" search-radius MSPE gain): row 19 is the base
#    value (1); rows 20-27 are gains relative to the previous row,
#    mirroring column K's relationship to the Area 3a block.
# ---------------------------------------------------------------------
$ws.Range("L19").Value2 = 1
for ($row = 20; $row -le 27; $row++) {
    $prev = $row - 1
    $ws.Range("L$row").Formula = "=(BD$prev-BD$row)/BD$row"
}

# ---------------------------------------------------------------------
# 6. Column width for AQ (col 43) now that it holds a header label.
# ---------------------------------------------------------------------
$ws.Columns.Item(43).ColumnWidth = 6.61

# ---------------------------------------------------------------------
# 7. Selection / view change.
# ---------------------------------------------------------------------
$ws.Range("N23").Select()
